# Add two student-score rows (id, name, score) below the existing header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 15211001 / 来吧快活吧 / 0
# Force column A to Text format first so the numeric-looking student id
# is stored as a string (shared string), not auto-converted to a number,
# then restore the cell's style so no stray formatting is left behind.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "15211001"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 2).Value = "来吧快活吧"
$ws.Cells.Item(2, 3).Value = 0

# Row 3: 15211002 / 你先都是你先 / 0
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "15211002"
$ws.Cells.Item(3, 1).Style = "Normal"
$ws.Cells.Item(3, 2).Value = "你先都是你先"
$ws.Cells.Item(3, 3).Value = 0
